# Weekly price-sheet update: a new week of "Melón" (Calameño / Tuna) price
# quotes is inserted at the top of the historical block (rows 294-299),
# pushing the existing rows 294-330 down to 300-336.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 294..330 down by 6 to make room for the new week's rows.
$ws.Range("294:299").Insert()

# Constant fields shared by every row in this market/category block.
$marketId = 11
$market   = "Vega Monumental Concepción"
$region   = "Bíobío"
$codreg   = 8
$catId    = 100112027
$category = "Melón"
$unit     = "`$/unidad"
$kgOrUnit = 1
$classif  = "Hortaliza"
$fecha    = 44946
$origen   = "Región de O'Higgins"

$newRows = @(
  @{ Row = 294; Variedad = "Calameño"; Calidad = "Primera"; Volumen = 1500; Min = 1000; Max = 1000; Prom = 1000 },
  @{ Row = 295; Variedad = "Calameño"; Calidad = "Segunda"; Volumen = 1000; Min =  800; Max =  800; Prom =  800 },
  @{ Row = 296; Variedad = "Calameño"; Calidad = "Tercera"; Volumen = 1000; Min =  600; Max =  600; Prom =  600 },
  @{ Row = 297; Variedad = "Tuna";     Calidad = "Primera"; Volumen = 1000; Min = 1000; Max = 1000; Prom = 1000 },
  @{ Row = 298; Variedad = "Tuna";     Calidad = "Segunda"; Volumen = 1000; Min =  800; Max =  800; Prom =  800 },
  @{ Row = 299; Variedad = "Tuna";     Calidad = "Tercera"; Volumen = 1000; Min =  600; Max =  600; Prom =  600 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $marketId
    $ws.Cells.Item($row, 2).Value  = $market
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $category
    $ws.Cells.Item($row, 8).Value  = $r.Variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Min
    $ws.Cells.Item($row, 12).Value = $r.Max
    $ws.Cells.Item($row, 13).Value = $r.Prom
    $ws.Cells.Item($row, 14).Value = $unit
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $kgOrUnit
    $ws.Cells.Item($row, 18).Value = $classif
}
